$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.362.82'
$ws.Range('E2').Value = '  +0.49%  '

$ws.Range('D3').Value = '3.266.23'
$ws.Range('E3').Value = '  +2.89%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '614.75'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.35%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.70'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.97%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = '3.264.12'
$ws.Range('E8').Value = '  +2.90%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.544'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.19%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.161'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.12%  '

$ws.Range('E11').Value = '  +1.51%  '

$ws.Range('E12').Value = '  -4.29%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.57%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '39.03'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.07%  '

$ws.Range('D15').Value = '3.805.34'
$ws.Range('E15').Value = '  +2.98%  '

$ws.Range('D16').Value = '66.425.47'
$ws.Range('E16').Value = '  +0.46%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.42'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.66%  '

$ws.Range('D18').Value = '3.269.86'
$ws.Range('E18').Value = '  +2.91%  '

$ws.Range('E19').Value = '  +1.38%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '505.21'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.67%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.53'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.23%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.755'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.69%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.14'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.74%  '

$ws.Range('E24').Value = '  -0.76%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '87.03'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.15%  '

$ws.Range('E26').Value = '  +0.02%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.05'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.76%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.22'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.78%  '

$ws.Range('E29').Value = '  +0.87%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.129'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +47.52%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.05'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.01%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.87'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.13%  '

$ws.Range('E33').Value = '  +0.33%  '

$ws.Range('E34').Value = '  -0.13%  '

$ws.Range('E35').Value = '  -2.52%  '

$ws.Range('E36').Value = '  -0.33%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.41'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +20.06%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '55.63'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.76%  '

$ws.Range('D39').Value = '0.0₃0790'
$ws.Range('E39').Value = '  +15.51%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '494.55'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.40%  '

$ws.Range('E41').Value = '  +1.64%  '

$ws.Range('E42').Value = '  +0.47%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.82'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.67%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.53'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.71%  '

$ws.Range('D45').Value = '3.014.79'
$ws.Range('E45').Value = '  +6.49%  '

$ws.Range('E46').Value = '  -1.24%  '

$ws.Range('E47').Value = '  +3.36%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.49'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.97%  '

$ws.Range('E49').Value = '  +2.37%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '121.26'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.84%  '
